$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the curated dimension/value cells (row 2-4, columns E and F)
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "iaest-measure:sexo"

$ws.Range("F3").Value = "medida"

$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("F4").Value = "xsd:int"

# Remove the now-unused mapping-file row entirely
$ws.Rows("5:5").Delete()
